$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

# Row 2
Set-TextValue 2 4 '59.195.27'
Set-TextValue 2 5 '  +4.92%  '

# Row 3
Set-TextValue 3 4 '3.322.46'
Set-TextValue 3 5 '  +2.20%  '

# Row 4
Set-TextValue 4 5 '  +0.08%  '

# Row 5
Set-TextValue 5 4 '408.72'
Set-TextValue 5 5 '  +2.44%  '

# Row 6
Set-TextValue 6 4 '112.94'
Set-TextValue 6 5 '  +1.94%  '

# Row 7
Set-TextValue 7 4 '0.586'
Set-TextValue 7 5 '  +4.96%  '

# Row 8
Set-TextValue 8 5 '  +0.06%  '

# Row 9
Set-TextValue 9 4 '0.634'
Set-TextValue 9 5 '  +2.64%  '

# Row 10
Set-TextValue 10 4 '40.14'
Set-TextValue 10 5 '  +1.82%  '

# Row 11
Set-TextValue 11 4 '0.0985'
Set-TextValue 11 5 '  +4.45%  '

# Row 13
Set-TextValue 13 4 '3.864.44'
Set-TextValue 13 5 '  +2.62%  '

# Row 14
Set-TextValue 14 4 '8.50'
Set-TextValue 14 5 '  +5.01%  '

# Row 15
Set-TextValue 15 4 '19.46'
Set-TextValue 15 5 '  +1.28%  '

# Row 16
Set-TextValue 16 4 '3.357.86'
Set-TextValue 16 5 '  +3.43%  '

# Row 17
Set-TextValue 17 5 '  -0.14%  '

# Row 18
Set-TextValue 18 4 '59.047.12'
Set-TextValue 18 5 '  +4.71%  '

# Row 19
Set-TextValue 19 4 '10.73'
Set-TextValue 19 5 '  -2.07%  '

# Row 20
Set-TextValue 20 4 '3.36'
Set-TextValue 20 5 '  +1.17%  '

# Row 21
Set-TextValue 21 4 '0.0000110'
Set-TextValue 21 5 '  +6.75%  '

# Row 22
Set-TextValue 22 4 '13.21'
Set-TextValue 22 5 '  +1.42%  '

# Row 23
Set-TextValue 23 4 '304.02'
Set-TextValue 23 5 '  +1.60%  '

# Row 24
Set-TextValue 24 4 '75.50'
Set-TextValue 24 5 '  +0.79%  '

# Row 25
Set-TextValue 25 4 '3.18'
Set-TextValue 25 5 '  -0.41%  '

# Row 26
Set-TextValue 26 2 'Kaspa'
Set-TextValue 26 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue 26 4 '0.183'
Set-TextValue 26 5 '  +8.04%  '

# Row 27
Set-TextValue 27 4 '28.52'
Set-TextValue 27 5 '  +1.40%  '

# Row 28
Set-TextValue 28 2 'LEO'
Set-TextValue 28 3 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue 28 4 '4.47'
Set-TextValue 28 5 '  +2.42%  '

# Row 29
Set-TextValue 29 4 '7.91'
Set-TextValue 29 5 '  -2.00%  '

# Row 30
Set-TextValue 30 4 '7.61'
Set-TextValue 30 5 '  +4.02%  '

# Row 31
Set-TextValue 31 5 '  +4.07%  '

# Row 32
Set-TextValue 32 4 '0.998'
Set-TextValue 32 5 '  -0.10%  '

# Row 33
Set-TextValue 33 4 '11.52'
Set-TextValue 33 5 '  +4.53%  '

# Row 34
Set-TextValue 34 4 '39.92'
Set-TextValue 34 5 '  +3.34%  '

# Row 35
Set-TextValue 35 4 '0.0511'
Set-TextValue 35 5 '  +5.06%  '

# Row 36
Set-TextValue 36 4 '52.14'
Set-TextValue 36 5 '  +1.27%  '

# Row 37
Set-TextValue 37 4 '2.10'
Set-TextValue 37 5 '  -1.87%  '

# Row 38
Set-TextValue 38 4 '3.13'
Set-TextValue 38 5 '  -0.26%  '

# Row 39
Set-TextValue 39 5 '  +0.02%  '

# Row 40
Set-TextValue 40 4 '3.40'
Set-TextValue 40 5 '  -3.29%  '

# Row 41
Set-TextValue 41 4 '138.18'
Set-TextValue 41 5 '  +3.37%  '

# Row 42
Set-TextValue 42 5 '  +2.36%  '

# Row 43
Set-TextValue 43 4 '1.92'
Set-TextValue 43 5 '  -0.89%  '

# Row 44
Set-TextValue 44 4 '17.00'
Set-TextValue 44 5 '  -2.88%  '

# Row 45
Set-TextValue 45 4 '3.93'
Set-TextValue 45 5 '  -0.84%  '

# Row 46
Set-TextValue 46 4 '0.283'
Set-TextValue 46 5 '  -0.35%  '

# Row 47
Set-TextValue 47 2 'WEMIXToken'
Set-TextValue 47 3 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 47 4 '2.27'
Set-TextValue 47 5 '  +8.52%  '

# Row 48
Set-TextValue 48 2 'EnergySwap'
Set-TextValue 48 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 48 4 '22.48'
Set-TextValue 48 5 '  +1.66%  '

# Row 49
Set-TextValue 49 4 '2.219.22'
Set-TextValue 49 5 '  +3.42%  '

# Row 50
Set-TextValue 50 5 '  -0.47%  '

# Row 51
Set-TextValue 51 5 '  -5.25%  '
